# Add release/1.0.2 to meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "release/1.0.2"
$ws.Range("B3:D3").Value = "X"
